$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).ClearContents()

$ws.Cells.Item(106, 8).Value = 2524.0527
$ws.Cells.Item(106, 9).Value = 2009.9375
$ws.Cells.Item(106, 10).Value = 5266
$ws.Cells.Item(106, 11).Value = 2009.9375
$ws.Cells.Item(106, 12).Value = 5266
$ws.Cells.Item(106, 13).Value = -1378.9375
$ws.Cells.Item(106, 14).Value = -6528

$ws.Cells.Item(112, 8).Value = 1345.7778
$ws.Cells.Item(112, 10).Value = 1345.7778
$ws.Cells.Item(112, 12).Value = 4037.3334
$ws.Cells.Item(112, 14).Value = -6253.3334

$ws.Cells.Item(121, 8).Value = 602.61536
$ws.Cells.Item(121, 10).Value = 602.61536
$ws.Cells.Item(121, 12).Value = 1807.84608
$ws.Cells.Item(121, 14).Value = -5301.84608

$ws.Cells.Item(130, 8).Value = 41884.285
$ws.Cells.Item(130, 10).Value = 41884.285
$ws.Cells.Item(130, 12).Value = 41884.285
$ws.Cells.Item(130, 14).Value = -51924.285

$ws.Cells.Item(137, 8).Value = 520285.28
$ws.Cells.Item(137, 10).Value = 2767.9636
$ws.Cells.Item(137, 12).Value = 8303.890800000001
$ws.Cells.Item(137, 14).Value = -13403.8908

$ws.Cells.Item(138, 8).Value = 3261.147
$ws.Cells.Item(138, 9).Value = 1727.2106
$ws.Cells.Item(138, 10).Value = 3855.9387
$ws.Cells.Item(138, 11).Value = 5181.6318
$ws.Cells.Item(138, 12).Value = 11567.8161
$ws.Cells.Item(138, 13).Value = -41.63180000000011
$ws.Cells.Item(138, 14).Value = -21847.8161

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3829.987
$ws.Cells.Item(32, 9).Value = 3883.1091
$ws.Cells.Item(32, 10).Value = 3697.182
$ws.Cells.Item(32, 11).Value = 3883.1091
$ws.Cells.Item(32, 12).Value = 3697.182
$ws.Cells.Item(32, 13).Value = -3596.1091
$ws.Cells.Item(32, 14).Value = -4271.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 4750
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 4750
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 4750
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).Value = -7746

$ws.Cells.Item(134, 8).Value = 3897.361
$ws.Cells.Item(134, 9).Value = 1394.1428
$ws.Cells.Item(134, 10).Value = 7401.8667
$ws.Cells.Item(134, 11).Value = 4182.428400000001
$ws.Cells.Item(134, 12).Value = 22205.6001
$ws.Cells.Item(134, 13).Value = -1647.428400000001
$ws.Cells.Item(134, 14).Value = -27275.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2061.6377
$ws.Cells.Item(31, 9).Value = 840.1613
$ws.Cells.Item(31, 10).Value = 3058.1052
$ws.Cells.Item(31, 11).Value = 840.1613
$ws.Cells.Item(31, 12).Value = 3058.1052
$ws.Cells.Item(31, 13).Value = -545.1613
$ws.Cells.Item(31, 14).Value = -3648.1052

$ws.Cells.Item(34, 8).Value = 2061.6377
$ws.Cells.Item(34, 9).Value = 840.1613
$ws.Cells.Item(34, 10).Value = 3058.1052
$ws.Cells.Item(34, 11).Value = 840.1613
$ws.Cells.Item(34, 12).Value = 3058.1052
$ws.Cells.Item(34, 13).Value = -638.1613
$ws.Cells.Item(34, 14).Value = -3462.1052

$ws.Cells.Item(58, 8).Value = 2376.878
$ws.Cells.Item(58, 9).Value = 1347.6061
$ws.Cells.Item(58, 10).Value = 6622.625
$ws.Cells.Item(58, 11).Value = 1347.6061
$ws.Cells.Item(58, 12).Value = 6622.625
$ws.Cells.Item(58, 13).Value = -1144.6061
$ws.Cells.Item(58, 14).Value = -7028.625

$ws.Cells.Item(99, 8).Value = 15388662
$ws.Cells.Item(99, 9).Value = 28573088
$ws.Cells.Item(99, 10).Value = 6833.3335
$ws.Cells.Item(99, 11).Value = 28573088
$ws.Cells.Item(99, 12).Value = 6833.3335
$ws.Cells.Item(99, 13).Value = -28571590
$ws.Cells.Item(99, 14).Value = -9829.333500000001

$ws.Cells.Item(126, 8).Value = 15388662
$ws.Cells.Item(126, 9).Value = 28573088
$ws.Cells.Item(126, 10).Value = 6833.3335
$ws.Cells.Item(126, 11).Value = 85719264
$ws.Cells.Item(126, 12).Value = 20500.0005
$ws.Cells.Item(126, 13).Value = -85716794
$ws.Cells.Item(126, 14).Value = -25440.0005

$ws.Cells.Item(132, 8).Value = 3524.2593
$ws.Cells.Item(132, 9).Value = 2839.875
$ws.Cells.Item(132, 10).Value = 8999.333000000001
$ws.Cells.Item(132, 11).Value = 8519.625
$ws.Cells.Item(132, 12).Value = 26997.999
$ws.Cells.Item(132, 13).Value = -5989.625
$ws.Cells.Item(132, 14).Value = -32057.999

$ws.Cells.Item(136, 8).Value = 2376.878
$ws.Cells.Item(136, 9).Value = 1347.6061
$ws.Cells.Item(136, 10).Value = 6622.625
$ws.Cells.Item(136, 11).Value = 4042.8183
$ws.Cells.Item(136, 12).Value = 19867.875
$ws.Cells.Item(136, 13).Value = -1492.8183
$ws.Cells.Item(136, 14).Value = -24967.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1406.7097
$ws.Cells.Item(68, 9).Value = 905.13043
$ws.Cells.Item(68, 10).Value = 1702.5128
$ws.Cells.Item(68, 11).Value = 2715.39129
$ws.Cells.Item(68, 12).Value = 5107.538399999999
$ws.Cells.Item(68, 13).Value = -1904.39129
$ws.Cells.Item(68, 14).Value = -6729.538399999999

$ws.Cells.Item(71, 8).Value = 1406.7097
$ws.Cells.Item(71, 9).Value = 905.13043
$ws.Cells.Item(71, 10).Value = 1702.5128
$ws.Cells.Item(71, 11).Value = 8146.173870000001
$ws.Cells.Item(71, 12).Value = 15322.6152
$ws.Cells.Item(71, 13).Value = -4090.173870000001
$ws.Cells.Item(71, 14).Value = -23434.6152

$ws.Cells.Item(107, 8).Value = 6681202.5
$ws.Cells.Item(107, 9).Value = 505.9091
$ws.Cells.Item(107, 10).Value = 11930321
$ws.Cells.Item(107, 11).Value = 1517.7273
$ws.Cells.Item(107, 12).Value = 35790963
$ws.Cells.Item(107, 13).Value = 402.2727
$ws.Cells.Item(107, 14).Value = -35794803

$ws.Cells.Item(131, 8).Value = 874.8
$ws.Cells.Item(131, 10).Value = 947.30334
$ws.Cells.Item(131, 12).Value = 2841.91002
$ws.Cells.Item(131, 14).Value = -12921.91002

$ws.Cells.Item(138, 8).Value = 2682.8572
$ws.Cells.Item(138, 9).Value = 926.6667
$ws.Cells.Item(138, 11).Value = 2780.0001
$ws.Cells.Item(138, 13).Value = 2359.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 9107815
$ws.Cells.Item(11, 9).Value = 19375012
$ws.Cells.Item(11, 10).Value = 1640762.9
$ws.Cells.Item(11, 11).Value = 19375012
$ws.Cells.Item(11, 12).Value = 1640762.9
$ws.Cells.Item(11, 13).Value = -19374873
$ws.Cells.Item(11, 14).Value = -1641040.9

$ws.Cells.Item(132, 8).Value = 4034.923
$ws.Cells.Item(132, 9).Value = 2319.182
$ws.Cells.Item(132, 10).Value = 5293.1333
$ws.Cells.Item(132, 11).Value = 6957.545999999999
$ws.Cells.Item(132, 12).Value = 15879.3999
$ws.Cells.Item(132, 13).Value = -4427.545999999999
$ws.Cells.Item(132, 14).Value = -20939.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2354.5
$ws.Cells.Item(100, 9).Value = 1950
$ws.Cells.Item(100, 11).Value = 1950
$ws.Cells.Item(100, 13).Value = -1409

$ws.Cells.Item(132, 8).Value = 3145.204
$ws.Cells.Item(132, 9).Value = 2413.4595
$ws.Cells.Item(132, 10).Value = 5401.4165
$ws.Cells.Item(132, 11).Value = 7240.3785
$ws.Cells.Item(132, 12).Value = 16204.2495
$ws.Cells.Item(132, 13).Value = -4710.3785
$ws.Cells.Item(132, 14).Value = -21264.2495

$ws.Cells.Item(136, 8).Value = 3339.1904
$ws.Cells.Item(136, 9).Value = 1544.4
$ws.Cells.Item(136, 10).Value = 5978.5884
$ws.Cells.Item(136, 11).Value = 4633.200000000001
$ws.Cells.Item(136, 12).Value = 17935.7652
$ws.Cells.Item(136, 13).Value = -2083.200000000001
$ws.Cells.Item(136, 14).Value = -23035.7652
